$wb = $excel.ActiveWorkbook

# --- "lines" sheet: insert a new "I_lim_A" column after v_nom_kv ---
$ws = $wb.Worksheets.Item("lines")

# Shift existing columns C:I (length_km .. is_pu) one column to the right,
# into D:J, working from the rightmost column first so values aren't clobbered.
for ($col = 9; $col -ge 3; $col--) {
    $srcHeaderVal = $ws.Cells.Item(1, $col).Value2
    $srcValueVal  = $ws.Cells.Item(2, $col).Value2
    $ws.Cells.Item(1, $col + 1).Value = $srcHeaderVal
    $ws.Cells.Item(2, $col + 1).Value = $srcValueVal
}

# New column C: I_lim_A
$ws.Cells.Item(1, 3).Value = "I_lim_A"
$ws.Cells.Item(2, 3).Value = 1000

# Move the "trafos" sheet selection off its previous cell first, then
# select "lines" last so "lines" ends up as the active tab.
$trafos = $wb.Worksheets.Item("trafos")
$trafos.Range("H6").Select()

$ws.Range("C1:C2").Select()
